# This workbook has two sheets: "Add-User" (sheet1) and "Upload-File" (sheet2).
# The edit reworks the Add-User test data: the header row is relabeled with
# human-readable captions, the single remaining data row now carries a new
# user's details, and the second data row (the old "C-TC002" case) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add-User")

# --- Header row (row 1) -----------------------------------------------
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "User Role"
$ws.Range("D1").Value = "Employee Name"
$ws.Range("E1").Value = "Status"
$ws.Range("F1").Value = "Password"
$ws.Range("G1").Value = "Confirm Password"

# --- Data row (row 2) - now represents test case C-TC001 --------------
$ws.Range("B2").Value = "timothy.amiano"
$ws.Range("C2").Value = "Admin"
$ws.Range("D2").Value = "Timothy Lewis Amiano"
$ws.Range("E2").Value = "Enabled"
$ws.Range("F2").Value = "password1"
$ws.Range("G2").Value = "password1"

# --- Remove the old second data row (C-TC002) --------------------------
$ws.Rows.Item(3).Delete()

# --- Column D is now wider to fit the longer employee name -------------
$ws.Columns.Item(4).ColumnWidth = 18.59

# --- Restore the selection left by the author after editing ------------
[void]$ws.Range("H19").Select()
